# Insert a new data row at row 192 (pushing existing rows 192:289 down to 193:290)
# and populate the new row with a fresh Tomate price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("192:192").Insert()

$ws.Range("A192").Value = 7
$ws.Range("B192").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C192").Value = 'Ñuble'
$ws.Range("D192").Value = 44466
$ws.Range("E192").Value = 16
$ws.Range("F192").Value = 100112020
$ws.Range("G192").Value = 'Tomate'
$ws.Range("H192").Value = 'Larga vida'
$ws.Range("I192").Value = 'Primera'
$ws.Range("J192").Value = 300
$ws.Range("K192").Value = 7500
$ws.Range("L192").Value = 8000
$ws.Range("M192").Value = 7750
$ws.Range("N192").Value = '$/caja 10 kilos'
$ws.Range("O192").Value = 'Región de Arica y Parinacota'
$ws.Range("P192").Value = 775
$ws.Range("Q192").Value = 10
$ws.Range("R192").Value = 'Hortaliza'
